# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 566 (shifting the existing
# rows 566-646 down to 567-647), keeping all other data untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 566; everything below shifts down one row.
$ws.Rows("566:566").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A566").Value = 3
$ws.Range("B566").Value = 'Femacal de La Calera'
$ws.Range("C566").Value = 'Coquimbo'
$ws.Range("D566").Value = 45077
$ws.Range("D566").NumberFormat = $ws.Range("D567").NumberFormat
$ws.Range("E566").Value = 5
$ws.Range("F566").Value = 100112017
$ws.Range("G566").Value = 'Apio'
$ws.Range("H566").Value = 'Americana (o)'
$ws.Range("I566").Value = 'Primera'
$ws.Range("J566").Value = 230
$ws.Range("K566").Value = 9000
$ws.Range("L566").Value = 9500
$ws.Range("M566").Value = 9261
$ws.Range("N566").Value = '$/docena de matas'
$ws.Range("O566").Value = 'Provincia de Limarí'
$ws.Range("P566").Value = 1544
$ws.Range("Q566").Value = 6
$ws.Range("R566").Value = 'Hortaliza'
